$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
Write-Host "theme name before:" $t.Name
try {
$t.Name = "Office Theme"
Write-Host "theme name after:" $t.Name
} catch {
Write-Host "ERR" $_.Exception.Message
}
